# Commit: "Commit in relation to performing of the dispatcher and performer process"
#
# On the "Constants" sheet, two new configuration entries are inserted:
#   1. PathTempDirectory   -> inserted as a new row right after "PathCustomerNameList"
#   2. TempFile_FileName   -> inserted as a new row right after
#                              "TemplateFile_WorksheetReviewSheetTemplate"
# Both insertions push every row below them down by one (full row insert),
# matching the dispatcher/performer hand-off described in the commit message
# (temp file + temp directory used to pass data between the two processes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# --- Insert #1: PathTempDirectory, right after PathCustomerNameList (row 21) ---
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).RowHeight = 14.25

$ws.Cells.Item(22, 1).Value = "PathTempDirectory"
$ws.Cells.Item(22, 2).Value = "C:\Users\martin.martinez\Documents\UiPath\temp"
$ws.Cells.Item(22, 3).Value = "path where the dispatcher process saved the temporary file with the information"

# --- Insert #2: TempFile_FileName, right after TemplateFile_WorksheetReviewSheetTemplate (now row 25) ---
$ws.Rows.Item(26).Insert()
$ws.Rows.Item(26).RowHeight = 14.25

$ws.Cells.Item(26, 1).Value = "TempFile_FileName"
$ws.Cells.Item(26, 2).Value = "temp.xlsx"
$ws.Cells.Item(26, 3).Value = "name of the excel where the dispatcher process saved the necessary information"

# --- Cosmetic: move the active selection to roughly where the author left it ---
$ws.Activate()
[void]$ws.Range("B18").Select()
